$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 326.6
$ws.Range("I11").Value = 326.6
$ws.Range("K11").Value = 326.6
$ws.Range("M11").Value = -186.6
$ws.Range("H39").Value = 326.2857
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6592
$ws.Range("H74").Value = 5298.3076
$ws.Range("I74").Value = 2959.3333
$ws.Range("K74").Value = 2959.3333
$ws.Range("M74").Value = -2023.3333
$ws.Range("H77").Value = 5298.3076
$ws.Range("I77").Value = 2959.3333
$ws.Range("K77").Value = 14796.6665
$ws.Range("M77").Value = -10116.6665
$ws.Range("H98").Value = 1267.8
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H122").Value = 1267.8
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("H131").Value = 4995.375
$ws.Range("I131").Value = 3389.6365
$ws.Range("J131").Value = 6354.077
$ws.Range("K131").Value = 10168.9095
$ws.Range("L131").Value = 19062.231
$ws.Range("M131").Value = -5128.9095
$ws.Range("N131").Value = -29142.231
$ws.Range("H135").Value = 848.25
$ws.Range("I135").Value = 751.73334
$ws.Range("J135").Value = 2296
$ws.Range("K135").Value = 6765.60006
$ws.Range("L135").Value = 20664
$ws.Range("M135").Value = -4230.60006
$ws.Range("N135").Value = -25734
$ws.Range("H137").Value = 3196.5356
$ws.Range("I137").Value = 3001.85
$ws.Range("K137").Value = 9005.549999999999
$ws.Range("M137").Value = -6455.549999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7932.2974
$ws.Range("I32").Value = 7160.4116
$ws.Range("K32").Value = 7160.4116
$ws.Range("M32").Value = -6873.4116
$ws.Range("H61").Value = 26320252
$ws.Range("I61").Value = 31253612
$ws.Range("J61").Value = 9000
$ws.Range("K61").Value = 31253612
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = -31253400
$ws.Range("N61").Value = -9424
$ws.Range("H136").Value = 26320252
$ws.Range("I136").Value = 31253612
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 93760836
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -93758286
$ws.Range("N136").Value = -32100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 150000
$ws.Range("J61").Value = 150000
$ws.Range("L61").Value = 150000
$ws.Range("N61").Value = -150626
$ws.Range("H107").Value = 3839.889
$ws.Range("I107").Value = 3768.3333
$ws.Range("K107").Value = 3768.3333
$ws.Range("M107").Value = -1848.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 72505.5
$ws.Range("I105").Value = 45000
$ws.Range("K105").Value = 45000
$ws.Range("M105").Value = -43253
$ws.Range("H107").Value = 2342.5715
$ws.Range("I107").Value = 1746.25
$ws.Range("K107").Value = 1746.25
$ws.Range("M107").Value = 173.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 412.5
$ws.Range("J29").Value = 516.6667
$ws.Range("L29").Value = 1550.0001
$ws.Range("N29").Value = -2104.0001
$ws.Range("H33").Value = 199.66667
$ws.Range("I33").Value = 50
$ws.Range("K33").Value = 300
$ws.Range("M33").Value = -17
$ws.Range("H88").Value = 4000
$ws.Range("J88").Value = 4000
$ws.Range("L88").Value = 12000
$ws.Range("N88").Value = -12856
$ws.Range("H91").Value = 4000
$ws.Range("J91").Value = 4000
$ws.Range("L91").Value = 12000
$ws.Range("N91").Value = -14964
$ws.Range("H92").Value = 457.7
$ws.Range("J92").Value = 457.7
$ws.Range("L92").Value = 1373.1
$ws.Range("N92").Value = -3869.1
$ws.Range("H131").Value = 28861.342
$ws.Range("J131").Value = 4435.2256
$ws.Range("L131").Value = 13305.6768
$ws.Range("N131").Value = -23385.6768
$ws.Range("H134").Value = 2777.9
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws.Range("H140").Value = 4672.25
$ws.Range("J140").Value = 4000
$ws.Range("L140").Value = 12000
$ws.Range("N140").Value = -22360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3977.5
$ws.Range("J80").Value = 3503
$ws.Range("L80").Value = 3503
$ws.Range("N80").Value = -5499
$ws.Range("H83").Value = 3977.5
$ws.Range("J83").Value = 3503
$ws.Range("L83").Value = 17515
$ws.Range("N83").Value = -27499
$ws.Range("H107").Value = 470.41177
$ws.Range("I107").Value = 413.8
$ws.Range("J107").Value = 895
$ws.Range("K107").Value = 413.8
$ws.Range("L107").Value = 895
$ws.Range("M107").Value = 1506.2
$ws.Range("N107").Value = -4735
$ws.Range("H113").Value = 3681.0667
$ws.Range("I113").Value = 3237.7407
$ws.Range("J113").Value = 4346.0557
$ws.Range("K113").Value = 3237.7407
$ws.Range("L113").Value = 4346.0557
$ws.Range("M113").Value = -1067.7407
$ws.Range("N113").Value = -8686.055700000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4250.2
$ws.Range("I22").Value = 1499.5
$ws.Range("K22").Value = 1499.5
$ws.Range("M22").Value = -1204.5
$ws.Range("H27").Value = 4250.2
$ws.Range("I27").Value = 1499.5
$ws.Range("K27").Value = 1499.5
$ws.Range("M27").Value = -1392.5
$ws.Range("H40").Value = 3094.182
$ws.Range("I40").Value = 3094.182
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3094.182
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2958.182
$ws.Range("N40").Value = $null
$ws.Range("H50").Value = 36000
$ws.Range("J50").Value = 36000
$ws.Range("L50").Value = 36000
$ws.Range("N50").Value = -37274
$ws.Range("H54").Value = 33998.668
$ws.Range("J54").Value = 33249.125
$ws.Range("L54").Value = 33249.125
$ws.Range("N54").Value = -34537.125
$ws.Range("H55").Value = 943.1
$ws.Range("I55").Value = 254.83333
$ws.Range("K55").Value = 254.83333
$ws.Range("M55").Value = -81.83332999999999
$ws.Range("H93").Value = 3053.5
$ws.Range("I93").Value = 2658
$ws.Range("K93").Value = 2658
$ws.Range("M93").Value = -1410
$ws.Range("H133").Value = 69447.5
$ws.Range("J133").Value = 69447.5
$ws.Range("L133").Value = 69447.5
$ws.Range("N133").Value = -74507.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 200014350
$ws.Range("I41").Value = 500009500
$ws.Range("J41").Value = 17599.666
$ws.Range("K41").Value = 500009500
$ws.Range("L41").Value = 17599.666
$ws.Range("M41").Value = -500009110
$ws.Range("N41").Value = -18379.666
$ws.Range("H81").Value = 1715.2858
$ws.Range("J81").Value = 1834.5
$ws.Range("L81").Value = 3669
$ws.Range("N81").Value = -5791
$ws.Range("H84").Value = 1715.2858
$ws.Range("J84").Value = 1834.5
$ws.Range("L84").Value = 18345
$ws.Range("N84").Value = -28953
$ws.Range("H122").Value = 113966.445
$ws.Range("I122").Value = 127587.25
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 382761.75
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -380311.75
$ws.Range("N122").Value = -19900
